# Fruta / hortaliza, semanal
# Insert a new data row at row 1006 (pushing existing rows 1006-1078 down to 1007-1079).
# The new row duplicates the data that currently sits in the last row (1078) of the
# sheet, except with an updated date (Fecha), matching the latest weekly price
# observation that was appended to the source dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 1078
$newRow = 1006

# Insert a blank row at position 1006; rows 1006..1078 shift down to 1007..1079.
$ws.Rows.Item($newRow).Insert()

# After the shift, the row that used to be 1078 is now 1079 - copy its values into
# the freshly inserted row 1006 (columns A through R).
$lastColCount = 18
for ($col = 1; $col -le $lastColCount; $col++) {
    $srcCell = $ws.Cells.Item($lastRow + 1, $col)
    $dstCell = $ws.Cells.Item($newRow, $col)
    $dstCell.Value2 = $srcCell.Value2
}

# Update the date (column D = 4) of the new row to the new observation date.
$ws.Cells.Item($newRow, 4).Value2 = 45021
